$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 163, shifting existing
# rows 163..256 down to 164..257 (dimension grows from R256 to R257).
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new weekly price record.
$ws.Cells.Item(163, 1).Value = 5
$ws.Cells.Item(163, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(163, 3).Value = "Maule"
$ws.Cells.Item(163, 4).Value = 44719
$ws.Cells.Item(163, 5).Value = 7
$ws.Cells.Item(163, 6).Value = 100112009
$ws.Cells.Item(163, 7).Value = "Acelga"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 500
$ws.Cells.Item(163, 11).Value = 2500
$ws.Cells.Item(163, 12).Value = 2500
$ws.Cells.Item(163, 13).Value = 2500
$ws.Cells.Item(163, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 625
$ws.Cells.Item(163, 17).Value = 4
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(163, 4).NumberFormat = $ws.Cells.Item(164, 4).NumberFormat
